# Generate Report for Handoff
# A brand-new handoff round has kicked off: the source markdown file was
# regenerated under a new GUID (d3e98c78-... -> c4a32591-...), producing a
# new xliff digest (5844b2e5... -> bba7337c...) and fresh handoff timestamps.
# Because this is a *new* handoff, the previously recorded target/handback
# info for each locale is no longer valid and is cleared out.

$wb = $excel.ActiveWorkbook

$oldGuid = "d3e98c78-52c1-48c9-870d-bc0ad3910f2d"
$newGuid = "c4a32591-c44e-4744-b89f-023e16e20042"

$oldDigest = "5844b2e5a7e0a2a2313eaf3fd2d3614df64b38a7"
$newDigest = "bba7337cdbcf1044003d7ecadcda38804bb592b3"

$newFileName   = "$newGuid.md"
$newPathName   = "e2e\$newGuid.md"

# ---------------------------------------------------------------------
# Overview sheet: file name / path / latest HO xliff generate date
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = $newFileName
$wsOverview.Range("B2").Value = $newPathName
$wsOverview.Range("G2").Value = "2016-09-03 03:03:59"

foreach ($hl in $wsOverview.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq "`$B`$2") {
        $hl.TextToDisplay = $newPathName
    }
}

# ---------------------------------------------------------------------
# Per-locale sheets (zh-cn, de-de): new handoff file / datetime, and the
# previous target + handback bookkeeping is reset for the new round.
# ---------------------------------------------------------------------
$locales = @{
    "zh-cn" = @{ HandoffDate = "2016-09-03 03:03:54" }
    "de-de" = @{ HandoffDate = "2016-09-03 03:03:59" }
}

foreach ($localeName in $locales.Keys) {
    $ws = $wb.Worksheets.Item($localeName)
    $info = $locales[$localeName]

    $newHandoffFile = "$newGuid.$newDigest.$localeName.xlf"

    # A: Source File Name (also hyperlinked)
    $ws.Range("A2").Value = $newFileName
    foreach ($hl in $ws.Hyperlinks) {
        $addr = $hl.Range.Address()
        if ($addr -eq "`$A`$2") {
            $hl.TextToDisplay = $newFileName
        }
    }

    # G/H: Latest Handoff File / Latest Handoff Datetime
    $ws.Range("G2").Value = $newHandoffFile
    $ws.Range("H2").Value = $info.HandoffDate

    # I2 previously hyperlinked to the old handoff markdown file - remove
    # the now-stale "Latest Target File" hyperlink and value entirely.
    foreach ($hl in $ws.Hyperlinks) {
        $addr = $hl.Range.Address()
        if ($addr -eq "`$I`$2") {
            $hl.Delete()
        }
    }
    $ws.Range("I2").Style = "Normal"
    $ws.Range("I2").Value = ""

    # J: Latest Handback File - cleared, no handback yet for this round
    $ws.Range("J2").Value = ""

    # K: Latest Handback DateTime - reset to the zero/epoch sentinel
    $ws.Range("K2").Value = "0001-01-01 00:00:00"

    # I/J column widths shrink now that they hold short/empty values
    # instead of long file names.
    $ws.Columns.Item(9).ColumnWidth = 17.833333333333336
    $ws.Columns.Item(10).ColumnWidth = 20.833333333333336
}
